# Regression data update for calorimetry test workbook (ds.2.dsc/data.xlsx)
#
# Summary of the change (from the OOXML diff):
#  - The component previously labelled "MF3H" is renamed to "T3H" everywhere.
#  - constants_evaluated gains a "Component" column, keeps only 1 data row,
#    and the Constant / St.Deviation values become text (kept numerically
#    identical) while Validity becomes "-Inf".
#  - enthalpies_calculated keeps only 1 data row (for "Comp"), and the
#    enthalpy value's sign flips from negative to positive.
#  - input_stoich_coefficients keeps only 1 data row, with stoichiometric
#    coefficients stored as text "1"/"1" instead of numbers.
#  - input_enthalpies clears the "PLP" data row (row 2) to empty cells.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Global rename: "MF3H" -> "T3H"
# ---------------------------------------------------------------------
$stoich = $wb.Worksheets.Item("input_stoich_coefficients")
$stoich.Cells.Item(1, 2).Value = "T3H"

$inputEnthalpies = $wb.Worksheets.Item("input_enthalpies")
$inputEnthalpies.Cells.Item(3, 1).Value = "T3H"

$setup = $wb.Worksheets.Item("setup")
$setup.Cells.Item(2, 2).Value = "T3H"

# ---------------------------------------------------------------------
# 2) constants_evaluated: add "Component" column, collapse to 1 data row,
#    store Constant/St.Deviation as text, Validity becomes "-Inf".
# ---------------------------------------------------------------------
$constEval = $wb.Worksheets.Item("constants_evaluated")

# Delete the now-superfluous rows 4 and 3 (bottom-up so indices stay valid).
$constEval.Rows.Item(4).Delete()
$constEval.Rows.Item(3).Delete()

# Shift existing headers from A:C to B:D, and add the new "Component" header.
$constEval.Cells.Item(1, 4).Value = "Validity"
$constEval.Cells.Item(1, 3).Value = "St.Deviation"
$constEval.Cells.Item(1, 2).Value = "Constant"
$constEval.Cells.Item(1, 1).Value = "Component"

# Row 2 values -- B2/C2 keep their numeric text but become strings.
$constEval.Cells.Item(2, 1).Value = "Comp"

$constEval.Cells.Item(2, 2).NumberFormat = "@"
$constEval.Cells.Item(2, 2).Value = "5.4061110496521"

$constEval.Cells.Item(2, 3).NumberFormat = "@"
$constEval.Cells.Item(2, 3).Value = "0.0205130354010432"

$constEval.Cells.Item(2, 4).Value = "-Inf"

# ---------------------------------------------------------------------
# 3) enthalpies_calculated: collapse to 1 data row ("Comp"), flip the
#    enthalpy value's sign from negative to positive.
# ---------------------------------------------------------------------
$enthCalc = $wb.Worksheets.Item("enthalpies_calculated")

$enthCalc.Rows.Item(4).Delete()
$enthCalc.Rows.Item(3).Delete()

$enthCalc.Cells.Item(2, 1).Value = "Comp"
$enthCalc.Cells.Item(2, 2).Value = 50.7744283493441
$enthCalc.Cells.Item(2, 3).Value = 0.534436612020058

# ---------------------------------------------------------------------
# 4) input_stoich_coefficients: collapse to 1 data row, coefficients
#    stored as text "1" instead of numbers.
# ---------------------------------------------------------------------
$stoich.Rows.Item(4).Delete()
$stoich.Rows.Item(3).Delete()

$stoich.Cells.Item(2, 1).NumberFormat = "@"
$stoich.Cells.Item(2, 1).Value = "1"

$stoich.Cells.Item(2, 2).NumberFormat = "@"
$stoich.Cells.Item(2, 2).Value = "1"

$stoich.Cells.Item(2, 3).Value = "Comp"

# ---------------------------------------------------------------------
# 5) input_enthalpies: clear the "PLP" row (row 2) to empty cells, but
#    keep the (now-blank) cells present -- touch a no-op formatting
#    property so the row/cells are not pruned entirely on save.
# ---------------------------------------------------------------------
$inputEnthalpies.Cells.Item(2, 1).Value = ""
$inputEnthalpies.Cells.Item(2, 2).Value = ""
$inputEnthalpies.Cells.Item(2, 1).Font.Bold = $false
$inputEnthalpies.Cells.Item(2, 2).Font.Bold = $false
